$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) cells per the
# refreshed cryptos snapshot. Column D values are forced to Text so
# numeric-looking strings (e.g. "1.004") are not auto-converted to
# numbers, matching the original inline-string cell contents.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.755.11"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +0.42%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.702.92"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +0.22%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  +0.30%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "317.11"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -0.38%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3940"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -0.64%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4048"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +0.25%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.523"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -1.62%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +0.29%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "53.65"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -1.37%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08893"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +0.71%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.536"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +3.14%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.72"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +1.20%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.166"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +6.70%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00001327"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -0.41%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.709.75"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +0.35%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "99.88"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -1.53%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07050"
$cell.Style = $origStyle

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.75"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.30%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.095"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +2.86%  "

$ws.Range("E22").Value = "  +0.21%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.46"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +2.01%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.753.81"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +0.47%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.249"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +5.71%  "

$ws.Range("E26").Value = "  +1.70%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.78"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +1.31%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "162.46"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("E29").Value = "  +16.00%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "136.34"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +1.29%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.179"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -1.34%  "

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.873"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +4.39%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08924"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +3.60%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.084"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -3.38%  "

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.985"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +1.71%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.12"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -4.20%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2763"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +0.18%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.62"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -0.91%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02797"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -0.13%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09189"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +1.13%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.465"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -0.06%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7734"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -0.75%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "16.02"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +2.10%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7214"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -0.92%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.573"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +1.48%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.221"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -0.40%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  +0.34%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.348"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -2.14%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "140.92"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -0.67%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "90.98"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +2.05%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07984"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -0.94%  "
